# Applies:
#  1. Removes the Heading2 paragraph style (i.e. drops the <w:pPr> block
#     entirely) from the "Introduction", "Amazon's Impact on Human
#     Behavior", "Amazon's Leadership Approach" and "Conclusion"
#     paragraphs.
#  2. Replaces the inline author-citations with generated reference ids:
#       (Alimahomed-Wilson and Reese) -> (Ref-f352447)
#       (Miller and Miller)  [1st]    -> (Ref-f836201)
#       (Zehndorfer)                  -> (Ref-f836201)
#       (Miller and Miller)  [2nd]    -> (Ref-f038989)
#       (Aguinis et al.)              -> (Ref-f038989)

$d = $word.ActiveDocument

function Strip-HeadingStyle($paraIndex, $expectedText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $runText = $r.Text
    if (-not $runText.StartsWith($expectedText)) {
        throw "Paragraph $paraIndex text mismatch: expected '$expectedText', got '$runText'"
    }
    # Escape any XML-sensitive characters in the run text (none expected
    # for these short headings, but keep this safe/generic).
    $escaped = $runText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    # Trim the trailing paragraph mark character captured by Range.Text.
    $escaped = $escaped.TrimEnd([char]13, [char]7)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r><w:t>' + $escaped + '</w:t></w:r></w:p></w:body>' +
           '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml) | Out-Null
}

# --- 1. Drop the Heading2 style from the four section headings ---------
Strip-HeadingStyle 2 "Introduction"
Strip-HeadingStyle 6 "Amazon's Impact on Human Behavior"
Strip-HeadingStyle 10 "Amazon's Leadership Approach"
Strip-HeadingStyle 16 "Conclusion"

# --- 2. Swap inline citations for generated reference ids --------------

function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $ok = $d.Paragraphs.Item($paraIndex).Range.Find.Execute(
        $findText, $true, $false, $false, $false, $false,
        $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Paragraph $paraIndex : could not find '$findText' to replace"
    }
}

# Paragraph 8: "An example that vividly illustrates Amazon's influence..."
Replace-InParagraph 8 "(Alimahomed-Wilson and Reese)" "(Ref-f352447)"

# Paragraph 12: "In examining Amazon's leadership style..."
Replace-InParagraph 12 "(Miller and Miller)" "(Ref-f836201)"
Replace-InParagraph 12 "(Zehndorfer)" "(Ref-f836201)"

# Paragraph 14: "Additionally, Amazon's focus on an ego-driven leadership..."
Replace-InParagraph 14 "(Miller and Miller)" "(Ref-f038989)"
Replace-InParagraph 14 "(Aguinis et al.)" "(Ref-f038989)"
